$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 30000
$ws.Range("J13").Value = 30000
$ws.Range("L13").Value = 30000
$ws.Range("N13").Value = -30338
$ws.Range("H76").Value = 3438.3333
$ws.Range("I76").Value = 3326.6667
$ws.Range("K76").Value = 3326.6667
$ws.Range("M76").Value = -3011.6667
$ws.Range("H79").Value = 3438.3333
$ws.Range("I79").Value = 3326.6667
$ws.Range("K79").Value = 3326.6667
$ws.Range("M79").Value = -2234.6667
$ws.Range("H112").Value = 1349.1305
$ws.Range("J112").Value = 1349.1305
$ws.Range("L112").Value = 4047.3915
$ws.Range("N112").Value = -6263.3915
$ws.Range("H129").Value = 1542.3265
$ws.Range("J129").Value = 1574.1489
$ws.Range("L129").Value = 4722.4467
$ws.Range("N129").Value = -14722.4467
$ws.Range("H132").Value = 349540.3
$ws.Range("I132").Value = 146318.3
$ws.Range("J132").Value = 771616.9
$ws.Range("K132").Value = 438954.9
$ws.Range("L132").Value = 2314850.7
$ws.Range("M132").Value = -436424.9
$ws.Range("N132").Value = -2319910.7
$ws.Range("H137").Value = 493497.34
$ws.Range("I137").Value = 1289567.8
$ws.Range("J137").Value = 2587.2334
$ws.Range("K137").Value = 3868703.4
$ws.Range("L137").Value = 7761.7002
$ws.Range("M137").Value = -3866153.4
$ws.Range("N137").Value = -12861.7002
$ws.Range("H138").Value = 3276.7258
$ws.Range("I138").Value = 1789.7778
$ws.Range("J138").Value = 3885.0227
$ws.Range("K138").Value = 5369.3334
$ws.Range("L138").Value = 11655.0681
$ws.Range("M138").Value = -229.3334000000004
$ws.Range("N138").Value = -21935.0681
$ws.Range("H141").Value = 6854.952
$ws.Range("I141").Value = 7052.3687
$ws.Range("J141").Value = 4979.5
$ws.Range("K141").Value = 21157.1061
$ws.Range("L141").Value = 14938.5
$ws.Range("M141").Value = -15977.1061
$ws.Range("N141").Value = -25298.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3888.4246
$ws.Range("I32").Value = 4005.2263
$ws.Range("J32").Value = 3578.9
$ws.Range("K32").Value = 4005.2263
$ws.Range("L32").Value = 3578.9
$ws.Range("M32").Value = -3718.2263
$ws.Range("N32").Value = -4152.9
$ws.Range("H102").Value = 1958.5714
$ws.Range("I102").Value = 1958.5714
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1958.5714
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -336.5714
$ws.Range("H122").Value = 2751.32
$ws.Range("I122").Value = 2387.3096
$ws.Range("K122").Value = 7161.9288
$ws.Range("M122").Value = -4711.9288
$ws.Range("H137").Value = 44591.5
$ws.Range("J137").Value = 44591.5
$ws.Range("L137").Value = 44591.5
$ws.Range("N137").Value = -54791.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H137").Value = 48770
$ws.Range("J137").Value = 48770
$ws.Range("L137").Value = 48770
$ws.Range("N137").Value = -58970
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2165.6462
$ws.Range("I31").Value = 894.5
$ws.Range("J31").Value = 3127.5945
$ws.Range("K31").Value = 894.5
$ws.Range("L31").Value = 3127.5945
$ws.Range("M31").Value = -599.5
$ws.Range("N31").Value = -3717.5945
$ws.Range("H34").Value = 2165.6462
$ws.Range("I34").Value = 894.5
$ws.Range("J34").Value = 3127.5945
$ws.Range("K34").Value = 894.5
$ws.Range("L34").Value = 3127.5945
$ws.Range("M34").Value = -692.5
$ws.Range("N34").Value = -3531.5945
$ws.Range("H58").Value = 2626.0286
$ws.Range("I58").Value = 1480.3448
$ws.Range("J58").Value = 8163.5
$ws.Range("K58").Value = 1480.3448
$ws.Range("L58").Value = 8163.5
$ws.Range("M58").Value = -1277.3448
$ws.Range("N58").Value = -8569.5
$ws.Range("H136").Value = 2626.0286
$ws.Range("I136").Value = 1480.3448
$ws.Range("J136").Value = 8163.5
$ws.Range("K136").Value = 4441.0344
$ws.Range("L136").Value = 24490.5
$ws.Range("M136").Value = -1891.0344
$ws.Range("N136").Value = -29590.5
$ws.Range("H139").Value = 38945
$ws.Range("J139").Value = 38945
$ws.Range("L139").Value = 38945
$ws.Range("N139").Value = -49225
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4600
$ws.Range("I63").Value = 4100
$ws.Range("J63").Value = 4877.778
$ws.Range("K63").Value = 12300
$ws.Range("L63").Value = 14633.334
$ws.Range("M63").Value = -11551
$ws.Range("N63").Value = -16131.334
$ws.Range("H66").Value = 4600
$ws.Range("I66").Value = 4100
$ws.Range("J66").Value = 4877.778
$ws.Range("K66").Value = 36900
$ws.Range("L66").Value = 43900.002
$ws.Range("M66").Value = -33156
$ws.Range("N66").Value = -51388.002
$ws.Range("H68").Value = 1407.6061
$ws.Range("I68").Value = 983.76
$ws.Range("J68").Value = 1666.0488
$ws.Range("K68").Value = 2951.28
$ws.Range("L68").Value = 4998.1464
$ws.Range("M68").Value = -2140.28
$ws.Range("N68").Value = -6620.1464
$ws.Range("H71").Value = 1407.6061
$ws.Range("I71").Value = 983.76
$ws.Range("J71").Value = 1666.0488
$ws.Range("K71").Value = 8853.84
$ws.Range("L71").Value = 14994.4392
$ws.Range("M71").Value = -4797.84
$ws.Range("N71").Value = -23106.4392
$ws.Range("H95").Value = 11999.75
$ws.Range("J95").Value = 11999.75
$ws.Range("L95").Value = 35999.25
$ws.Range("N95").Value = -40117.25
$ws.Range("H107").Value = 6681207.5
$ws.Range("I107").Value = 513.90625
$ws.Range("J107").Value = 11652887
$ws.Range("K107").Value = 1541.71875
$ws.Range("L107").Value = 34958661
$ws.Range("M107").Value = 378.28125
$ws.Range("N107").Value = -34962501
$ws.Range("H113").Value = 5435403
$ws.Range("I113").Value = 657.3333
$ws.Range("J113").Value = 11364217
$ws.Range("K113").Value = 1971.9999
$ws.Range("L113").Value = 34092651
$ws.Range("M113").Value = 198.0001
$ws.Range("N113").Value = -34096991
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 19286.732
$ws.Range("I43").Value = 1939.6
$ws.Range("J43").Value = 27960.3
$ws.Range("K43").Value = 1939.6
$ws.Range("L43").Value = 27960.3
$ws.Range("M43").Value = -1788.6
$ws.Range("N43").Value = -28262.3
$ws.Range("H46").Value = 24522.428
$ws.Range("J46").Value = 25755.334
$ws.Range("L46").Value = 25755.334
$ws.Range("N46").Value = -26067.334
$ws.Range("H137").Value = 31868.334
$ws.Range("J137").Value = 43736.668
$ws.Range("L137").Value = 43736.668
$ws.Range("N137").Value = -53936.668
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2395.818
$ws.Range("I100").Value = 1960
$ws.Range("J100").Value = 2759
$ws.Range("K100").Value = 1960
$ws.Range("L100").Value = 2759
$ws.Range("M100").Value = -1419
$ws.Range("N100").Value = -3841
$ws.Range("H132").Value = 3432.0454
$ws.Range("I132").Value = 2639.3333
$ws.Range("J132").Value = 5810.1816
$ws.Range("K132").Value = 7917.999899999999
$ws.Range("L132").Value = 17430.5448
$ws.Range("M132").Value = -5387.999899999999
$ws.Range("N132").Value = -22490.5448
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 52750
$ws.Range("J130").Value = 52750
$ws.Range("L130").Value = 52750
$ws.Range("N130").Value = -62790
$ws.Range("H136").Value = 3819.5862
$ws.Range("I136").Value = 1671.6666
$ws.Range("J136").Value = 6120.9287
$ws.Range("K136").Value = 5014.9998
$ws.Range("L136").Value = 18362.7861
$ws.Range("M136").Value = -2464.9998
$ws.Range("N136").Value = -23462.7861
